# Fix mojibake "Â±" (U+00C2 U+00B1) -> "±" (U+00B1) introduced by a
# double UTF-8 encoding, in the plus-minus statistics columns
# (f1_score_weighted, training_time, test_time) for data rows 2-17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mojibake = [string][char]0x00C2 + [char]0x00B1
$fixed = [string][char]0x00B1

for ($row = 2; $row -le 17; $row++) {
    foreach ($col in @("B", "C", "D")) {
        $cell = $ws.Range("$col$row")
        $value = $cell.Value2
        if ($value -ne $null -and $value.ToString().Contains($mojibake)) {
            $cell.Value2 = $value.ToString().Replace($mojibake, $fixed)
        }
    }
}
